$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.140614986419678
$ws.Range("B1").Value = 5.94737434387207
$ws.Range("C1").Value = 2.488657236099243
$ws.Range("D1").Value = 1.144976496696472
$ws.Range("E1").Value = 0.8179686069488525
